# Generate Report for handoff
# Update the "Latest Handoff Datetime" (column D) for the row belonging to
# the "1f5afceb-ce4f-4ad6-8e89-3fffebbcefdb" file (row 5) on both the
# zh-cn and de-de localization-status sheets to reflect the new handoff.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-14 15:12:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-14 15:13:18"
